# Apply the scripted text replacements to the active document.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-04-11 Friday"; new = "2025-04-12 Saturday"},
    @{old = "130×2=260";         new = "684×8=5472"},
    @{old = "978×6=5868";        new = "149×2=298"},
    @{old = "282×8=2256";        new = "848×5=4240"},
    @{old = "501×5=2505";        new = "824×8=6592"},
    @{old = "340×3=1020";        new = "217×7=1519"},
    @{old = "429×5=2145";        new = "488×7=3416"},
    @{old = "854×7=5978";        new = "847×8=6776"},
    @{old = "988×2=1976";        new = "901×7=6307"},
    @{old = "884×6=5304";        new = "775×9=6975"},
    @{old = "676×5=3380";        new = "462×7=3234"},
    @{old = "952×7=6664";        new = "588×3=1764"},
    @{old = "420×9=3780";        new = "238×6=1428"},
    @{old = "959×5=4795";        new = "514×3=1542"},
    @{old = "764×7=5348";        new = "745×5=3725"},
    @{old = "414×5=2070";        new = "650×8=5200"},
    @{old = "574×9=5166";        new = "816×8=6528"},
    @{old = "527×2=1054";        new = "812×2=1624"},
    @{old = "599×8=4792";        new = "609×2=1218"},
    @{old = "526×2=1052";        new = "259×2=518"},
    @{old = "852×7=5964";        new = "746×2=1492"},
    @{old = "746×6=4476";        new = "129×3=387"},
    @{old = "509×9=4581";        new = "872×7=6104"},
    @{old = "705×9=6345";        new = "954×4=3816"},
    @{old = "829×2=1658";        new = "921×6=5526"},
    @{old = "151×7=1057";        new = "803×8=6424"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
